$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly price log (rows 2..63, row 1 = header).
# Two new daily records need to be inserted as the new rows 9 and 10,
# pushing the former rows 9..63 down to rows 11..65.

$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# New row 9: "Primera" quality record dated 2022-03-04
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value = "Maule"
$ws.Cells.Item(9, 4).Value = 44624
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100101
$ws.Cells.Item(9, 8).Value = "Berries"
$ws.Cells.Item(9, 9).Value = 100101001
$ws.Cells.Item(9, 10).Value = "Arándano (blue)"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 120
$ws.Cells.Item(9, 14).Value = 3300
$ws.Cells.Item(9, 15).Value = 3300
$ws.Cells.Item(9, 16).Value = 3300
$ws.Cells.Item(9, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Linares"
$ws.Cells.Item(9, 19).Value = 1650
$ws.Cells.Item(9, 20).Value = 2

# New row 10: "Segunda" quality record, same date 2022-03-04
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44624
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101001
$ws.Cells.Item(10, 10).Value = "Arándano (blue)"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 100
$ws.Cells.Item(10, 14).Value = 3000
$ws.Cells.Item(10, 15).Value = 3000
$ws.Cells.Item(10, 16).Value = 3000
$ws.Cells.Item(10, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Linares"
$ws.Cells.Item(10, 19).Value = 1500
$ws.Cells.Item(10, 20).Value = 2
